$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns for the data rows stay text-typed so that
# numeric-looking strings (e.g. "1.00", "0.0000226") are not coerced into numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '61.660.82'
$ws.Range("E2").Value = '  -2.29%  '
$ws.Range("D3").Value = '2.975.44'
$ws.Range("E3").Value = '  -2.67%  '
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").Value = '587.61'
$ws.Range("E5").Value = '  +0.63%  '
$ws.Range("D6").Value = '141.72'
$ws.Range("E6").Value = '  -6.77%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.06%  '
$ws.Range("D8").Value = '0.521'
$ws.Range("E8").Value = '  -2.68%  '
$ws.Range("D9").Value = '2.972.15'
$ws.Range("E9").Value = '  -2.78%  '
$ws.Range("E10").Value = '  -6.12%  '
$ws.Range("D11").Value = '5.74'
$ws.Range("E11").Value = '  -1.90%  '
$ws.Range("D12").Value = '0.458'
$ws.Range("E12").Value = '  +2.01%  '
$ws.Range("D13").Value = '0.0000226'
$ws.Range("E13").Value = '  -3.85%  '
$ws.Range("D14").Value = '33.97'
$ws.Range("E14").Value = '  -6.26%  '
$ws.Range("E15").Value = '  +1.30%  '
$ws.Range("D16").Value = '3.449.24'
$ws.Range("E16").Value = '  -3.07%  '
$ws.Range("D17").Value = '7.01'
$ws.Range("E17").Value = '  -2.00%  '
$ws.Range("D18").Value = '61.519.51'
$ws.Range("D19").Value = '2.953.94'
$ws.Range("E19").Value = '  -3.35%  '
$ws.Range("D20").Value = '449.92'
$ws.Range("E20").Value = '  -6.55%  '
$ws.Range("D21").Value = '13.86'
$ws.Range("E21").Value = '  -3.44%  '
$ws.Range("D22").Value = '0.682'
$ws.Range("E22").Value = '  -3.81%  '
$ws.Range("D23").Value = '7.31'
$ws.Range("E23").Value = '  -3.33%  '
$ws.Range("D24").Value = '81.06'
$ws.Range("E24").Value = '  -1.21%  '
$ws.Range("D25").Value = '12.11'
$ws.Range("E25").Value = '  -4.73%  '
$ws.Range("D26").Value = '2.15'
$ws.Range("E26").Value = '  -10.96%  '
$ws.Range("E27").Value = '  +0.08%  '
$ws.Range("D28").Value = '9.70'
$ws.Range("E28").Value = '  -8.29%  '
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("D30").Value = '2.62'
$ws.Range("E30").Value = '  -1.53%  '
$ws.Range("D31").Value = '6.87'
$ws.Range("E31").Value = '  -7.30%  '
$ws.Range("D32").Value = '2.06'
$ws.Range("E32").Value = '  -6.84%  '
$ws.Range("D33").Value = '27.19'
$ws.Range("E33").Value = '  -2.23%  '
$ws.Range("D34").Value = '0.107'
$ws.Range("E34").Value = '  -3.88%  '
$ws.Range("D35").Value = '1.01'
$ws.Range("E35").Value = '  -5.20%  '
$ws.Range("D36").Value = '0.0₃0773'
$ws.Range("E36").Value = '  -5.72%  '
$ws.Range("D37").Value = '5.67'
$ws.Range("E37").Value = '  -4.47%  '
$ws.Range("D38").Value = '2.08'
$ws.Range("E38").Value = '  -5.92%  '
$ws.Range("B39").Value = 'Cosmos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D39").Value = '9.14'
$ws.Range("E39").Value = '  -0.87%  '
$ws.Range("B40").Value = 'OKB'
$ws.Range("C40").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D40").Value = '50.10'
$ws.Range("E40").Value = '  -0.80%  '
$ws.Range("D41").Value = '0.119'
$ws.Range("E41").Value = '  +2.92%  '
$ws.Range("D42").Value = '2.80'
$ws.Range("E42").Value = '  -13.09%  '
$ws.Range("D43").Value = '387.25'
$ws.Range("E43").Value = '  -9.77%  '
$ws.Range("D44").Value = '0.0352'
$ws.Range("E44").Value = '  -3.01%  '
$ws.Range("D45").Value = '2.716.54'
$ws.Range("E45").Value = '  -4.59%  '
$ws.Range("D46").Value = '0.262'
$ws.Range("E46").Value = '  -8.76%  '
$ws.Range("D47").Value = '37.00'
$ws.Range("E47").Value = '  -2.93%  '
$ws.Range("D48").Value = '129.76'
$ws.Range("E48").Value = '  +1.92%  '
$ws.Range("D50").Value = '0.108'
$ws.Range("E50").Value = '  -1.83%  '
$ws.Range("D51").Value = '2.16'
$ws.Range("E51").Value = '  -2.18%  '

# Restore the original (default) cell style now that the text values are set.
$ws.Range("D2:E51").Style = "Normal"
